$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sprint2 sheet: reflect US09 completion + re-baseline the "Est Time" column
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Sprint2")

# US08 (row 3) Est Time 26 -> 25
$ws.Cells.Item(3, 6).Value = 25

# US09 (row 4) Est Time 27 -> 25, plus newly recorded Act Size / Act Time / Completed
$ws.Cells.Item(4, 6).Value = 25
$ws.Cells.Item(4, 7).Value = 62
$ws.Cells.Item(4, 8).Value = 35
$ws.Cells.Item(4, 9).Value = "Completed"

# Remaining stories re-baselined Est Time
$ws.Cells.Item(5, 6).Value = 25
$ws.Cells.Item(6, 6).Value = 35
$ws.Cells.Item(7, 6).Value = 25

# Selection / activation moves to H9 on Sprint2, and Sprint2 becomes the active tab
$ws.Range("H9").Select()
$excel.ActiveWindow.ActiveSheet.Name | Out-Null
$ws.Activate()

$wb.Windows.Item(1).ActiveSheet.Range("H9").Select() | Out-Null
